$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the "Rovers need to know the plateau size" paragraph had its
# sentence split into two runs around a stray "_GoBack" bookmark (an artifact
# of where the cursor was when the doc was last saved in Word). Re-Find &
# Replace the full sentence so it collapses back into a single run and the
# now-orphaned bookmark goes away.
# ---------------------------------------------------------------------------
$sentence = "When a rover receives its instructions, it will calculate its route and if it will collide or fall off the plateau it will refuse to move. "
$d.Content.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, $sentence, 2)

# ---------------------------------------------------------------------------
# Change 2: append a new "Controller" class block to the Class Diagram
# section, right after the RobotArm class' last method, mirroring the
# formatting used by the other class headings (bold + underlined name
# followed by plain member lines). A trailing empty paragraph carries the
# document's "_GoBack" last-edit bookmark, matching where Word would leave it
# after typing this content.
# ---------------------------------------------------------------------------
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "^\+move\(x: number, y: number, z: number\)\r?$") {
        $targetIndex = $i + 1
        break
    }
}

$headingPara = $d.Paragraphs($targetIndex)
$headingPara.Range.InsertAfter("Controller")

# Make room for the member lines straight away (before the heading gets its
# bold/underline styling) so the new paragraphs inherit plain formatting.
$d.Paragraphs($targetIndex).Range.InsertParagraphAfter()
$d.Paragraphs($targetIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs($targetIndex + 2).Range.InsertParagraphAfter()
$d.Paragraphs($targetIndex + 3).Range.InsertParagraphAfter()

$d.Paragraphs($targetIndex + 1).Range.InsertAfter("+readInput(String[])")
$d.Paragraphs($targetIndex + 2).Range.InsertAfter("+createRover(id:String):Rover")
$d.Paragraphs($targetIndex + 3).Range.InsertAfter("+sendCommand(Rover, String)")

# Style the "Controller" heading line (bold + single underline), matching
# e.g. the "RobotArm" heading above it.
$headingRange = $d.Paragraphs($targetIndex).Range
$headingRange.Font.Bold = $true
$headingRange.Font.Underline = 1
$headingRange.LanguageID = "en-US"

# The trailing blank paragraph keeps the document's "_GoBack" bookmark
# (Word re-creates this at the last editing location on save).
$lastBlankPara = $d.Paragraphs($targetIndex + 4)
$d.Bookmarks.Add("_GoBack", $lastBlankPara.Range)
